$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: update title (D26)
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 32: update title (D32) and link (E32)
$ws.Range("D32").Value = "Knowledge Distillation"
$ws.Range("E32").Value = "https://dodonam.tistory.com/364"

# Row 46: update title (D46) and link (E46)
$ws.Range("D46").Value = "폐렴(Pneumonia) 종류"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/459"
